# Regenerate save_data: update column G ("K") values per recomputed
# std/mean and s_vals. Commit message: "regen save_data to use K instead
# of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 1
    6  = 3
    7  = 3
    8  = 2
    9  = 2
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 2
    15 = 2
    16 = 0
    17 = 0
    18 = 2
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 4
    26 = 1
    27 = 1
    28 = 2
    29 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
